$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Fitness values (column C) for rows 2-12 as per diff
$ws.Range("C2").Value = 3875
$ws.Range("C3").Value = 3954
$ws.Range("C4").Value = 4121
$ws.Range("C5").Value = 4121
$ws.Range("C6").Value = 4196
$ws.Range("C7").Value = 4318
$ws.Range("C8").Value = 4374
$ws.Range("C9").Value = 4595
$ws.Range("C10").Value = 4707
$ws.Range("C11").Value = 4707
$ws.Range("C12").Value = 4840
